# Fruta / hortaliza, semanal
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dateFmt = "YYYY-MM-DD HH:MM:SS"

# This weekly refresh appends the two previous observations (the current
# contents of rows 5 and 6, about to be overwritten below) to the bottom
# of the table as rows 7 and 8, then adds a brand-new row 9 for this
# week's third observation, before finally updating rows 5 and 6 in place
# with the newest figures.

# Row 7 <- old row 5
$ws.Cells.Item(7, 1).Value = 12
$ws.Cells.Item(7, 2).Value = "Mapocho Venta Directa de Santiago"
$ws.Cells.Item(7, 3).Value = "Metropolitana"
$ws.Cells.Item(7, 4).Value = 44376
$ws.Cells.Item(7, 4).NumberFormat = $dateFmt
$ws.Cells.Item(7, 5).Value = 13
$ws.Cells.Item(7, 6).Value = 100112026
$ws.Cells.Item(7, 7).Value = "Haba"
$ws.Cells.Item(7, 8).Value = "Sin especificar"
$ws.Cells.Item(7, 9).Value = "Primera"
$ws.Cells.Item(7, 10).Value = 15
$ws.Cells.Item(7, 11).Value = 12000
$ws.Cells.Item(7, 12).Value = 12000
$ws.Cells.Item(7, 13).Value = 12000
$ws.Cells.Item(7, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(7, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(7, 16).Value = 480
$ws.Cells.Item(7, 17).Value = 25
$ws.Cells.Item(7, 18).Value = "Hortaliza"

# Row 8 <- old row 6
$ws.Cells.Item(8, 1).Value = 12
$ws.Cells.Item(8, 2).Value = "Mapocho Venta Directa de Santiago"
$ws.Cells.Item(8, 3).Value = "Metropolitana"
$ws.Cells.Item(8, 4).Value = 44418
$ws.Cells.Item(8, 4).NumberFormat = $dateFmt
$ws.Cells.Item(8, 5).Value = 13
$ws.Cells.Item(8, 6).Value = 100112026
$ws.Cells.Item(8, 7).Value = "Haba"
$ws.Cells.Item(8, 8).Value = "Sin especificar"
$ws.Cells.Item(8, 9).Value = "Primera"
$ws.Cells.Item(8, 10).Value = 12
$ws.Cells.Item(8, 11).Value = 15000
$ws.Cells.Item(8, 12).Value = 15000
$ws.Cells.Item(8, 13).Value = 15000
$ws.Cells.Item(8, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(8, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(8, 16).Value = 600
$ws.Cells.Item(8, 17).Value = 25
$ws.Cells.Item(8, 18).Value = "Hortaliza"

# Row 9 <- new observation
$ws.Cells.Item(9, 1).Value = 12
$ws.Cells.Item(9, 2).Value = "Mapocho Venta Directa de Santiago"
$ws.Cells.Item(9, 3).Value = "Metropolitana"
$ws.Cells.Item(9, 4).Value = 44432
$ws.Cells.Item(9, 4).NumberFormat = $dateFmt
$ws.Cells.Item(9, 5).Value = 13
$ws.Cells.Item(9, 6).Value = 100112026
$ws.Cells.Item(9, 7).Value = "Haba"
$ws.Cells.Item(9, 8).Value = "Sin especificar"
$ws.Cells.Item(9, 9).Value = "Primera"
$ws.Cells.Item(9, 10).Value = 15
$ws.Cells.Item(9, 11).Value = 14000
$ws.Cells.Item(9, 12).Value = 14000
$ws.Cells.Item(9, 13).Value = 14000
$ws.Cells.Item(9, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(9, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(9, 16).Value = 560
$ws.Cells.Item(9, 17).Value = 25
$ws.Cells.Item(9, 18).Value = "Hortaliza"

# Row 5: update to this week's figures.
$ws.Cells.Item(5, 4).Value = 44435
$ws.Cells.Item(5, 11).Value = 14000
$ws.Cells.Item(5, 12).Value = 14000
$ws.Cells.Item(5, 13).Value = 14000
$ws.Cells.Item(5, 16).Value = 560

# Row 6: update to this week's figures.
$ws.Cells.Item(6, 4).Value = 44435
$ws.Cells.Item(6, 10).Value = 15
$ws.Cells.Item(6, 11).Value = 14000
$ws.Cells.Item(6, 12).Value = 14000
$ws.Cells.Item(6, 13).Value = 14000
$ws.Cells.Item(6, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(6, 16).Value = 560
